$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: TS001TC002TSP001 - Enter STDCIF in the Function Id
# Previously recorded as FAIL (browser closed); now PASS with real output.
$ws.Range("L7").Value = "PASS"
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = "Filled: Function Id"
$ws.Range("O7").Value = "screenshots/STEP_6.png"
$ws.Range("P7").Value = "page_sources/STEP_6_source.html"

# Row 8: TS001TC002TSP002 - Click the Go button
$ws.Range("L8").Value = "PASS"
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = "Clicked: Go"
$ws.Range("O8").Value = "screenshots/STEP_7.png"
$ws.Range("P8").Value = "page_sources/STEP_7_source.html"

# Row 9: TS001TC002TSP003 - Click New
$ws.Range("L9").Value = "PASS"
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = "Clicked: New"
$ws.Range("O9").Value = "screenshots/STEP_8.png"
$ws.Range("P9").Value = "page_sources/STEP_8_source.html"

# Row 10: TS001TC002TSP004 - Click "P" to generate customer no
$ws.Range("L10").Value = "PASS"
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = "Clicked: P"
$ws.Range("O10").Value = "screenshots/STEP_9.png"
$ws.Range("P10").Value = "page_sources/STEP_9_source.html"

# Row 11: TS001TC002TSP005 - enter full name
$ws.Range("L11").Value = "PASS"
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = "Filled: Full Name"
$ws.Range("O11").Value = "screenshots/STEP_10.png"
$ws.Range("P11").Value = "page_sources/STEP_10_source.html"

# Row 12: TS001TC002TSP006 - enter short name
$ws.Range("L12").Value = "PASS"
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = "Filled: Short Name"
$ws.Range("O12").Value = "screenshots/STEP_11.png"
$ws.Range("P12").Value = "page_sources/STEP_11_source.html"

# Row 13: TS001TC002TSP007 - enter customer category
$ws.Range("L13").Value = "PASS"
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = "Filled: Customer Category"
$ws.Range("O13").Value = "screenshots/STEP_12.png"
$ws.Range("P13").Value = "page_sources/STEP_12_source.html"
